$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 4.27
$ws.Range("C3").Value = 4.27
$ws.Range("B5").Value = 4.28
$ws.Range("B8").Value = 4.29
$ws.Range("B9").Value = 4.29
$ws.Range("B11").Value = 4.29
